# "Moved broom to Transform"
#
# The slide shows the R-for-Data-Science workflow diagram (Import / Tidy /
# Transform / Visualise / Model / Communicate). The "broom" hex-sticker
# picture, originally sitting next to "readr" under the Import step, is
# relocated into the big rounded-rectangle "Transform" zone. A handful of
# neighbouring shapes shift slightly to make room, the rounded rectangle is
# sent to the back of the z-order (so the broom now draws on top of it
# instead of below), and the trailing empty run on four of the step labels
# is trimmed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Move the "broom" picture (was beside Import, top-left) into the
#        Transform zone in the middle of the diagram. -----------------------
$broom = $s.Shapes.Item("Picture 48")
$broom.Left = 354.18001
$broom.Top = 324.0396

# --- 2. Send the big rounded rectangle ("Transform" zone background) to the
#        back of the z-order, so it now sits behind the broom picture that
#        used to be drawn after it. msoSendToBack = 1. ---------------------
$rect = $s.Shapes.Item("Rounded Rectangle 10")
$rect.ZOrder(1)

# --- 3. Nudge the "Import" label up slightly and drop its now-redundant
#        trailing empty run. ------------------------------------------------
$importBox = $s.Shapes.Item("TextBox 3")
$importRange = $importBox.TextFrame.TextRange
$importRange.Delete()
$importRange.Text = "Import"
$importBox.Left = 86.2074
$importBox.Top = 26.3493

# --- 4. Drop the redundant trailing empty run on the other three step
#        labels too (Tidy, Visualise, Model) -- no position change for these.
$tidyRange = $s.Shapes.Item("TextBox 4").TextFrame.TextRange
$tidyRange.Delete()
$tidyRange.Text = "Tidy"

$visualiseRange = $s.Shapes.Item("TextBox 6").TextFrame.TextRange
$visualiseRange.Delete()
$visualiseRange.Text = "Visualise"

$modelRange = $s.Shapes.Item("TextBox 7").TextFrame.TextRange
$modelRange.Delete()
$modelRange.Text = "Model"

# --- 5. Shift the second broom-ish picture (bottom-left, "Picture 49") a
#        bit to the right. --------------------------------------------------
$pic49 = $s.Shapes.Item("Picture 49")
$pic49.Left = 69.309
$pic49.Top = 62.02662

# --- 6. Lengthen/move the Tidy-side up-arrow ("Up Arrow 79") downward. -----
$arrow79 = $s.Shapes.Item("Up Arrow 79")
$arrow79.Top = 211.1691

# --- 7. Reposition/resize the Model-side up-arrow ("Up Arrow 84"). ---------
$arrow84 = $s.Shapes.Item("Up Arrow 84")
$arrow84.Left = 530.29
$arrow84.Top = 334.3616
$arrow84.Width = 16.441
$arrow84.Height = 107.7166
